# Add season-record columns (Wins, Losses, Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells, styled like the existing header row (bold/bordered).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Season record is the same for every player on the roster.
$wins = 74
$losses = 88
$ties = 0

for ($r = 2; $r -le 41; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins
    $ws.Cells.Item($r, 31).Value = $losses
    $ws.Cells.Item($r, 32).Value = $ties
}
